$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.710.05"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.898.58"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.08"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4832"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2856"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06561"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.001.77"
$ws.Range("E10").Value = "  +8.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07478"
$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.77"
$ws.Range("E12").Value = "  +2.78%  "

$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.21"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6689"
$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.711.34"
$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.265.40"
$ws.Range("E17").Value = "  +8.04%  "

$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007612"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "231.32"
$ws.Range("E21").Value = "  +2.74%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.242"
$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.98"
$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.358"
$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.80"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.968"
$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.404"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1018"
$ws.Range("E30").Value = "  +10.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.358"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05123"
$ws.Range("E33").Value = "  +2.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.222"
$ws.Range("E34").Value = "  +6.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7603"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.711"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01886"
$ws.Range("E37").Value = "  +3.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.660"
$ws.Range("E38").Value = "  +1.96%  "

$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.39"
$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4312"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.720"
$ws.Range("E44").Value = "  -4.04%  "

$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.85"
$ws.Range("E46").Value = "  +0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1276"
$ws.Range("E47").Value = "  -3.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.495"
$ws.Range("E48").Value = "  -3.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.984"
$ws.Range("E49").Value = "  +2.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.97"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("E51").Value = "  +0.45%  "

